# Generate Report for Handoff
# Updates the "Latest Handoff Date" / "Latest Handoff Datetime" cells for the
# 92d1389b-129a-498c-905a-c88a62160c79.md file after a fresh handoff was generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 7 corresponds to 92d1389b-129a-498c-905a-c88a62160c79.md
# Column D = "Latest Handoff Date"
$wsOverview.Range("D7").Value = "2016-03-21 22:37:41"

# zh-cn sheet: row 7 corresponds to 92d1389b-129a-498c-905a-c88a62160c79.md
# Column E = "Latest Handoff Datetime"
$wsZhCn.Range("E7").Value = "2016-03-21 22:37:37"

# de-de sheet: row 7 corresponds to 92d1389b-129a-498c-905a-c88a62160c79.md
# Column E = "Latest Handoff Datetime"
$wsDeDe.Range("E7").Value = "2016-03-21 22:37:41"
